$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "67.620.46"
Set-TextValue "E2" "  -1.09%  "
Set-TextValue "D3" "3.780.41"
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "595.31"
Set-TextValue "E5" "  -0.01%  "
Set-TextValue "D6" "166.38"
Set-TextValue "E6" "  -0.29%  "
Set-TextValue "D7" "3.779.30"
Set-TextValue "E7" "  +0.88%  "
Set-TextValue "E8" "  +0.10%  "
Set-TextValue "D9" "0.519"
Set-TextValue "E9" "  +0.02%  "
Set-TextValue "E10" "  +0.25%  "
Set-TextValue "D11" "6.36"
Set-TextValue "E11" "  -1.88%  "
Set-TextValue "E12" "  +0.09%  "
Set-TextValue "D13" "0.0000255"
Set-TextValue "E13" "  -1.24%  "
Set-TextValue "D14" "36.24"
Set-TextValue "E14" "  +0.65%  "
Set-TextValue "D15" "4.414.37"
Set-TextValue "E15" "  +0.89%  "
Set-TextValue "D16" "3.765.91"
Set-TextValue "E16" "  +0.46%  "
Set-TextValue "D17" "18.44"
Set-TextValue "E17" "  +3.13%  "
Set-TextValue "D18" "67.584.97"
Set-TextValue "E18" "  -1.04%  "
Set-TextValue "E19" "  +0.19%  "
Set-TextValue "E20" "  -0.16%  "
Set-TextValue "D21" "10.09"
Set-TextValue "E21" "  -5.67%  "
Set-TextValue "D22" "456.54"
Set-TextValue "E22" "  -2.36%  "
Set-TextValue "E23" "  +0.07%  "
Set-TextValue "E24" "  +7.92%  "
Set-TextValue "D26" "11.91"
Set-TextValue "E26" "  -0.59%  "
Set-TextValue "E27" "  -2.17%  "
Set-TextValue "E28" "  -0.07%  "
Set-TextValue "D30" "2.79"
Set-TextValue "E30" "  +0.34%  "
Set-TextValue "D31" "7.27"
Set-TextValue "E31" "  -0.30%  "
Set-TextValue "D32" "29.77"
Set-TextValue "E32" "  -0.22%  "
Set-TextValue "E33" "  +0.88%  "
Set-TextValue "E34" "  +0.12%  "
Set-TextValue "E35" "  -0.18%  "
Set-TextValue "D36" "3.732.73"
Set-TextValue "E36" "  +0.79%  "
Set-TextValue "E37" "  -0.88%  "
Set-TextValue "E38" "  -0.86%  "
Set-TextValue "E39" "  -1.11%  "
Set-TextValue "E40" "  -0.61%  "
Set-TextValue "D41" "5.75"
Set-TextValue "E41" "  -0.73%  "
Set-TextValue "D42" "0.999"
Set-TextValue "E42" "  -0.02%  "
Set-TextValue "D44" "45.37"
Set-TextValue "E44" "  +5.46%  "
Set-TextValue "E45" "  -1.48%  "
Set-TextValue "D46" "47.11"
Set-TextValue "E46" "  +2.59%  "
Set-TextValue "E47" "  -2.69%  "
Set-TextValue "D48" "148.41"
Set-TextValue "E48" "  +1.20%  "
Set-TextValue "E49" "  -4.23%  "
Set-TextValue "D50" "389.55"
Set-TextValue "E50" "  -0.13%  "
Set-TextValue "D51" "25.61"
Set-TextValue "E51" "  +0.77%  "
